# "Add ui and ux" - refresh the members list:
#  - drop the old "ABCD" row and the old first "Abu Inshah" row
#  - keep/renumber the remaining three member rows (ram, ram, Abu Inshah)
#  - the "ram" row now shows the "Wealth Manager" designation, and the
#    trailing "Abu Inshah" row reverts to "Health insurance advisor"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two obsolete rows; Excel shifts rows 4-6 up into 2-4 and
# shrinks the used range/dimension to A1:G4 automatically.
$ws.Rows("2:3").Delete()

# Make sure every data cell in the refreshed range is stored as text
# (phone numbers, epoch ids and ISO timestamps would otherwise be
# reinterpreted as numbers) before writing the values.
$ws.Range("A2:G4").NumberFormat = "@"

$ws.Range("A2").Value = "ram"
$ws.Range("B2").Value = "7449085120"
$ws.Range("C2").Value = "aiautomationhig@gmail.com"
$ws.Range("D2").Value = "Wealth Manager"
$ws.Range("E2").Value = "uploads/ram_1752311161933.jpeg"
$ws.Range("F2").Value = "1752311162001"
$ws.Range("G2").Value = "2025-07-12T09:06:02.001Z"

$ws.Range("A3").Value = "ram"
$ws.Range("B3").Value = "7449085120"
$ws.Range("C3").Value = "selvasuresh460@gmail.com"
$ws.Range("D3").Value = "Health insurance advisor"
$ws.Range("E3").Value = "uploads/ram_1752314093239.jpeg"
$ws.Range("F3").Value = "1752314093256"
$ws.Range("G3").Value = "2025-07-12T09:54:53.256Z"

$ws.Range("A4").Value = "Abu Inshah"
$ws.Range("B4").Value = "7449085120"
$ws.Range("C4").Value = "wealthplusacademy@gmail.com"
$ws.Range("D4").Value = "Health insurance advisor"
$ws.Range("E4").Value = "uploads/abu_inshah_1752314719383.jpeg"
$ws.Range("F4").Value = "1752314719399"
$ws.Range("G4").Value = "2025-07-12T10:05:19.399Z"

# Drop the temporary "@" number format again so the cells fall back to the
# workbook's default (unstyled) look, matching the original formatting.
$ws.Range("A2:G4").Style = "Normal"
